$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "data as of" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 19:10"

# Update country statistics with the latest data pull
# (rows stay sorted by total cases descending; Argelia/Serbia and
#  Groenlandia/Islas Malvinas swap order as their values changed)
$ws.Range("B4").Value = 3119823
$ws.Range("C4").Value = 22739
$ws.Range("D4").Value = 1359753
$ws.Range("E4").Value = 1625853
$ws.Range("G4").Value = 245
$ws.Range("H4").Value = 134217
$ws.Range("B5").Value = 1683738
$ws.Range("C5").Value = 9083
$ws.Range("E5").Value = 498703
$ws.Range("G5").Value = 245
$ws.Range("H5").Value = 67113
$ws.Range("B6").Value = 766273
$ws.Range("C6").Value = 22792
$ws.Range("D6").Value = 475480
$ws.Range("E6").Value = 269660
$ws.Range("G6").Value = 480
$ws.Range("H6").Value = 21133
$ws.Range("B9").Value = 303083
$ws.Range("C9").Value = 2064
$ws.Range("D9").Value = 271703
$ws.Range("E9").Value = 24807
$ws.Range("G9").Value = 139
$ws.Range("H9").Value = 6573
$ws.Range("B14").Value = 242149
$ws.Range("C14").Value = 193
$ws.Range("D14").Value = 193640
$ws.Range("E14").Value = 13595
$ws.Range("G14").Value = 15
$ws.Range("H14").Value = 34914
$ws.Range("B18").Value = 208938
$ws.Range("C18").Value = 1041
$ws.Range("D18").Value = 187511
$ws.Range("E18").Value = 16145
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 5282
$ws.Range("B23").Value = 106366
$ws.Range("C23").Value = 199
$ws.Range("D23").Value = 70161
$ws.Range("E23").Value = 27472
$ws.Range("G23").Value = 22
$ws.Range("H23").Value = 8733
$ws.Range("B28").Value = 73858
$ws.Range("C28").Value = 302
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = 5482
$ws.Range("B48").Value = 33175
$ws.Range("C48").Value = 953
$ws.Range("D48").Value = 18315
$ws.Range("E48").Value = 14516
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 344
$ws.Range("A62").Value = "Argelia"
$ws.Range("B62").Value = 17348
$ws.Range("C62").Value = 469
$ws.Range("D62").Value = 12329
$ws.Range("E62").Value = 4041
$ws.Range("G62").Value = 10
$ws.Range("H62").Value = 978
$ws.Range("A63").Value = "Serbia"
$ws.Range("B63").Value = 17076
$ws.Range("C63").Value = 357
$ws.Range("D63").Value = 13366
$ws.Range("E63").Value = 3369
$ws.Range("G63").Value = 11
$ws.Range("H63").Value = 341
$ws.Range("B92").Value = 5459
$ws.Range("C92").Value = 281
$ws.Range("D92").Value = 2349
$ws.Range("E92").Value = 3088
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 22
$ws.Range("B96").Value = 4650
$ws.Range("C96").Value = 47
$ws.Range("E96").Value = 484
$ws.Range("B111").Value = 2358
$ws.Range("C111").Value = 10
$ws.Range("D111").Value = 1597
$ws.Range("E111").Value = 641
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 120
$ws.Range("B112").Value = 2094
$ws.Range("C112").Value = 13
$ws.Range("E112").Value = 116
$ws.Range("B129").Value = 1221
$ws.Range("C129").Value = 16
$ws.Range("D129").Value = 1050
$ws.Range("E129").Value = 121
$ws.Range("B138").Value = 1008
$ws.Range("C138").Value = 3
$ws.Range("E138").Value = 150
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
